# Correction in SA algorithm and 746 logs.
# Updates the recorded "Fitness" values (column C) for generations 0-94
# (rows 2-96) of run 17, per the corrected simulated-annealing results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C4").Value = 8445
$ws.Range("C5:C6").Value = 8372
$ws.Range("C7:C9").Value = 8010
$ws.Range("C10:C29").Value = 7892
$ws.Range("C30:C36").Value = 7870
$ws.Range("C37:C96").Value = 7293
